$wb = $excel.ActiveWorkbook

# =====================================================================
# Sheet "Prix Spot": append a new day column (BT) with header "24-aug"
# and its 24 hourly prices, mirroring the existing BS ("23-aug") column.
# =====================================================================
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the header style (bold font, thin border, centered alignment)
# from the previous header cell (BS1) onto the new header cell (BT1).
$wsSpot.Range("BS1").Copy()
$wsSpot.Range("BT1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsSpot.Range("BT1").Value = "24-aug"

$wsSpot.Range("BT2").Value = 100.5
$wsSpot.Range("BT3").Value = 95
$wsSpot.Range("BT4").Value = 88.98999999999999
$wsSpot.Range("BT5").Value = 85.66
$wsSpot.Range("BT6").Value = 83.61
$wsSpot.Range("BT7").Value = 86.53
$wsSpot.Range("BT8").Value = 90
$wsSpot.Range("BT9").Value = 82.51000000000001
$wsSpot.Range("BT10").Value = 67.06
$wsSpot.Range("BT11").Value = 4.56
$wsSpot.Range("BT12").Value = -0.01
$wsSpot.Range("BT13").Value = -0.07000000000000001
$wsSpot.Range("BT14").Value = -1.98
$wsSpot.Range("BT15").Value = -7.45
$wsSpot.Range("BT16").Value = -7.4
$wsSpot.Range("BT17").Value = -0.1
$wsSpot.Range("BT18").Value = 0
$wsSpot.Range("BT19").Value = 12.93
$wsSpot.Range("BT20").Value = 61.15
$wsSpot.Range("BT21").Value = 98.56999999999999
$wsSpot.Range("BT22").Value = 114.99
$wsSpot.Range("BT23").Value = 114.32
$wsSpot.Range("BT24").Value = 107.5
$wsSpot.Range("BT25").Value = 100.52

# =====================================================================
# Sheet "Gaz": append a new daily quote row (69) for 2025-08-22.
# =====================================================================
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date cell to be entered as literal text (matching every
# other date cell in column A) instead of being auto-parsed as a date
# serial number; ClearFormats() afterwards drops the temporary "Text"
# number format so the cell keeps the sheet's default style.
$wsGaz.Range("A69").NumberFormat = "@"
$wsGaz.Range("A69").Value = "2025-08-22"
$wsGaz.Range("A69").ClearFormats()

$wsGaz.Range("B69").Value = 32.2

# =====================================================================
# Sheet "CO2": append a new daily quote row (69) for 2025-08-22, with
# the price left blank (not yet published).
# =====================================================================
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A69").NumberFormat = "@"
$wsCO2.Range("A69").Value = "2025-08-22"
$wsCO2.Range("A69").ClearFormats()

# Touch B69 so the (empty) cell exists in the sheet, then clear it back
# out to a truly blank value.
$wsCO2.Range("B69").NumberFormat = "@"
$wsCO2.Range("B69").Value = "x"
$wsCO2.Range("B69").Value = ""
$wsCO2.Range("B69").ClearFormats()
